$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an explicit text value, preventing Excel from
# auto-converting numeric-looking strings (e.g. "1.000") into numbers,
# and without leaving a residual explicit cell style behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "30.069.33"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").Value = "1.920.36"
$ws.Range("E3").Value = "  +2.45%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.86%  "
Set-TextValue $ws.Range("D5") "329.76"
$ws.Range("E5").Value = "  +4.42%  "
Set-TextValue $ws.Range("D6") "0.9995"
$ws.Range("E6").Value = "  -0.87%  "
Set-TextValue $ws.Range("D7") "0.5247"
$ws.Range("E7").Value = "  +3.11%  "
Set-TextValue $ws.Range("D8") "0.4055"
$ws.Range("E8").Value = "  +3.88%  "
Set-TextValue $ws.Range("D9") "0.08485"
$ws.Range("E9").Value = "  +1.48%  "
Set-TextValue $ws.Range("D10") "1.129"
$ws.Range("E10").Value = "  +2.32%  "
Set-TextValue $ws.Range("D11") "42.90"
$ws.Range("E11").Value = "  +3.09%  "
$ws.Range("E12").Value = "  +9.03%  "
Set-TextValue $ws.Range("D13") "6.384"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").Value = "1.921.07"
$ws.Range("E14").Value = "  +2.47%  "
Set-TextValue $ws.Range("D15") "7.401"
$ws.Range("E15").Value = "  +1.87%  "
Set-TextValue $ws.Range("D16") "1.000"
$ws.Range("E16").Value = "  -0.94%  "
Set-TextValue $ws.Range("D17") "96.52"
$ws.Range("E17").Value = "  +5.69%  "
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  -0.25%  "
Set-TextValue $ws.Range("D20") "18.24"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("E21").Value = "  -0.69%  "
Set-TextValue $ws.Range("D22") "6.070"
$ws.Range("E22").Value = "  +2.42%  "
$ws.Range("D23").Value = "30.077.92"
$ws.Range("E23").Value = "  +5.50%  "
Set-TextValue $ws.Range("D24") "11.27"
$ws.Range("E24").Value = "  +1.48%  "
Set-TextValue $ws.Range("D25") "2.224"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "2.141.91"
$ws.Range("E26").Value = "  +2.54%  "
Set-TextValue $ws.Range("D27") "21.17"
$ws.Range("E27").Value = "  +2.59%  "
Set-TextValue $ws.Range("D28") "160.22"
$ws.Range("E28").Value = "  -1.09%  "
Set-TextValue $ws.Range("D29") "2.459"
$ws.Range("E29").Value = "  +3.01%  "
Set-TextValue $ws.Range("D30") "129.47"
$ws.Range("E30").Value = "  +2.93%  "
$ws.Range("E31").Value = "  +4.15%  "
Set-TextValue $ws.Range("D32") "0.1060"
$ws.Range("E32").Value = "  +1.39%  "
Set-TextValue $ws.Range("D33") "6.112"
$ws.Range("E33").Value = "  +5.97%  "
Set-TextValue $ws.Range("D34") "3.655"
$ws.Range("E34").Value = "  +1.26%  "
Set-TextValue $ws.Range("D35") "0.02521"
$ws.Range("E35").Value = "  +2.58%  "
Set-TextValue $ws.Range("D36") "0.06608"
$ws.Range("E36").Value = "  +1.07%  "
Set-TextValue $ws.Range("D38") "9.036"
$ws.Range("E38").Value = "  +2.53%  "
Set-TextValue $ws.Range("D39") "1.237"
Set-TextValue $ws.Range("D40") "5.217"
$ws.Range("E40").Value = "  +3.35%  "
Set-TextValue $ws.Range("D41") "0.6568"
Set-TextValue $ws.Range("D42") "11.74"
$ws.Range("E42").Value = "  +5.89%  "
Set-TextValue $ws.Range("D43") "1.244"
$ws.Range("E43").Value = "  -0.01%  "
Set-TextValue $ws.Range("D44") "0.6194"
$ws.Range("E44").Value = "  +3.09%  "
Set-TextValue $ws.Range("D45") "13.33"
$ws.Range("E45").Value = "  +2.41%  "
Set-TextValue $ws.Range("D46") "3.776"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("E47").Value = "  +3.60%  "
Set-TextValue $ws.Range("D48") "125.98"
$ws.Range("E48").Value = "  +3.41%  "
Set-TextValue $ws.Range("D49") "1.244"
$ws.Range("E49").Value = "  +2.43%  "
Set-TextValue $ws.Range("D50") "80.17"
$ws.Range("E50").Value = "  +4.98%  "
Set-TextValue $ws.Range("D51") "1.153"
$ws.Range("E51").Value = "  +0.85%  "
